# handhelds_mtbenson_market workbook update
# - Corrects ingredient / allergen / local-ingredient text for several items
# - Fills in the "Nutrition Label" placeholder values with real per-item labels
# - Widens / re-sizes a few table columns and switches the whole table to
#   left-aligned text (was center-aligned)
# - Nudges the sheet's scroll position so column C is the first visible column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Corrected cell text (Ingredients / Allergens / LocalIngredients / Nutrition Label)
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Tuna / Red Onion / Celery / Pickles / Dill / Lemon / Garlic / Salt&Pepper / Lettuce / Mayo / Multigrain Bread"

$ws.Range("D4").Value = "Fresh Start Tomato, Island City Bread"

$ws.Range("B6").Value = "Ham / Swiss Cheese / Red Onion / Pickles / Tomato / Lettuce / Honey Dijon Aioli / Multigrain Bread"
$ws.Range("D6").Value = "Fresh Start Tomato, Island City Bread"
$ws.Range("F6").Value = "Ham_&_Cheese"

$ws.Range("D7").Value = "Fresh Start Tomato, Island City Bread"
$ws.Range("F7").Value = "Roast_Beef"

$ws.Range("B8").Value = "Falafel / Red Onion / Cucumber / Tomato / Feta / Lettuce / Tzatziki / Sundried Tomato Tortilla"
$ws.Range("F8").Value = "Falafel_Wrap"

$ws.Range("B9").Value = "Salami / Swiss Cheese / Tomato / Onion / Lettuce / Garlic Aioli / Sourdough Club "
$ws.Range("C9").Value = "Wheat, gluten, milk, eggs, sulphites, mustard."
$ws.Range("D9").Value = "Fresh Start Tomato, Castle Cheese Swiss, Nesvog Meats Salami"
$ws.Range("F9").Value = "Italian_Sandwich"

# ---------------------------------------------------------------------------
# 2. Column width changes (B, C, D)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 85.16666666666667
$ws.Columns.Item(3).ColumnWidth = 13.833333333333334
$ws.Columns.Item(4).ColumnWidth = 24.5

# ---------------------------------------------------------------------------
# 3. Re-align the whole table from center to left (header + data rows)
# ---------------------------------------------------------------------------
$full = $ws.Range("A1:G9")
$full.HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# 4. Scroll the sheet so column C is the left-most visible column
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
